$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.114.17'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '3.508.70'
$ws.Range("E3").Value = '  -2.41%  '
$ws.Range("D4").Value = '''1.00'
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = '''201.22'
$ws.Range("E5").Value = '  +4.21%  '
$ws.Range("D6").Value = '''549.73'
$ws.Range("E6").Value = '  -5.05%  '
$ws.Range("D7").Value = '3.497.72'
$ws.Range("E7").Value = '  -2.65%  '
$ws.Range("D8").Value = '''0.604'
$ws.Range("E8").Value = '  -2.22%  '
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("D10").Value = '''0.656'
$ws.Range("E10").Value = '  -3.28%  '
$ws.Range("D11").Value = '''61.04'
$ws.Range("E11").Value = '  +11.88%  '
$ws.Range("D12").Value = '''0.144'
$ws.Range("E12").Value = '  -4.39%  '
$ws.Range("D13").Value = '''0.0000276'
$ws.Range("E13").Value = '  +0.38%  '
$ws.Range("D14").Value = '''9.81'
$ws.Range("E14").Value = '  -1.49%  '
$ws.Range("D15").Value = '4.063.83'
$ws.Range("E15").Value = '  -2.30%  '
$ws.Range("D16").Value = '3.495.20'
$ws.Range("E16").Value = '  -2.57%  '
$ws.Range("E17").Value = '  -1.09%  '
$ws.Range("D18").Value = '''18.43'
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("D19").Value = '66.785.97'
$ws.Range("E19").Value = '  -0.85%  '
$ws.Range("D20").Value = '''11.84'
$ws.Range("E20").Value = '  -4.05%  '
$ws.Range("E21").Value = '  -3.79%  '
$ws.Range("D22").Value = '''387.88'
$ws.Range("E22").Value = '  -3.68%  '
$ws.Range("B23").Value = 'RenderToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D23").Value = '''12.03'
$ws.Range("E23").Value = '  -9.86%  '
$ws.Range("B24").Value = 'PancakeSwap'
$ws.Range("C24").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D24").Value = '''4.02'
$ws.Range("E24").Value = '  -4.75%  '
$ws.Range("D25").Value = '''82.35'
$ws.Range("E25").Value = '  -4.03%  '
$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").Value = '''2.80'
$ws.Range("E26").Value = '  -4.22%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").Value = '''11.95'
$ws.Range("E27").Value = '  -5.22%  '
$ws.Range("B28").Value = 'Toncoin'
$ws.Range("C28").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D28").Value = '''3.70'
$ws.Range("E28").Value = '  -2.58%  '
$ws.Range("B29").Value = 'Filecoin'
$ws.Range("C29").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D29").Value = '''8.88'
$ws.Range("E29").Value = '  -3.24%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").Value = '''30.69'
$ws.Range("E30").Value = '  -2.16%  '
$ws.Range("B31").Value = 'NEARProtocol'
$ws.Range("C31").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D31").Value = '''7.33'
$ws.Range("E31").Value = '  -9.39%  '
$ws.Range("B32").Value = 'Bittensor'
$ws.Range("C32").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D32").Value = '''670.52'
$ws.Range("E32").Value = '  +0.22%  '
$ws.Range("B33").Value = 'Cosmos'
$ws.Range("C33").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D33").Value = '''11.71'
$ws.Range("E33").Value = '  -4.31%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").Value = '''63.48'
$ws.Range("E34").Value = '  -1.05%  '
$ws.Range("B35").Value = 'Hedera'
$ws.Range("C35").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D35").Value = '''0.110'
$ws.Range("E35").Value = '  -5.53%  '
$ws.Range("B36").Value = 'InjectiveProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D36").Value = '''39.55'
$ws.Range("E36").Value = '  -7.81%  '
$ws.Range("B37").Value = 'TheGraph'
$ws.Range("C37").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D37").Value = '''0.412'
$ws.Range("E37").Value = '  -3.42%  '
$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").Value = '''1.00'
$ws.Range("E38").Value = '  +0.03%  '
$ws.Range("B39").Value = 'ThetaToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D39").Value = '''3.06'
$ws.Range("E39").Value = '  -1.74%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D40").Value = '3.112.34'
$ws.Range("E40").Value = '  -1.36%  '
$ws.Range("B41").Value = 'FirstDigitalUSD'
$ws.Range("C41").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D41").Value = '''0.997'
$ws.Range("E41").Value = '  -0.16%  '
$ws.Range("B42").Value = 'Kaspa'
$ws.Range("C42").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D42").Value = '''0.129'
$ws.Range("E42").Value = '  -3.71%  '
$ws.Range("B43").Value = 'PEPE'
$ws.Range("C43").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D43").Value = '0.0₃0707'
$ws.Range("E43").Value = '  -10.43%  '
$ws.Range("B44").Value = 'dogwifhat'
$ws.Range("C44").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D44").Value = '''2.82'
$ws.Range("E44").Value = '  +15.95%  '
$ws.Range("D45").Value = '''2.55'
$ws.Range("E45").Value = '  -13.60%  '
$ws.Range("B46").Value = 'WEMIXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D46").Value = '''2.68'
$ws.Range("E46").Value = '  +6.96%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").Value = '''0.0398'
$ws.Range("E47").Value = '  -4.88%  '
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D48").Value = '''0.127'
$ws.Range("E48").Value = '  -3.35%  '
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").Value = '''2.97'
$ws.Range("E49").Value = '  -5.09%  '
$ws.Range("D50").Value = '''136.58'
$ws.Range("E50").Value = '  -4.22%  '
$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").Value = '''8.29'
$ws.Range("E51").Value = '  -5.80%  '
